# Applies the changes described by the commit:
#  - Corrected two mistyped "Secteur" names in the Meteo sheet
#    (shared strings used by column B / "Secteur"):
#      "Torra di Murtella" -> "A Torra di Murtella"
#      "Maffalcu"           -> "Malfalcu"
#  - Widened column B (Secteur) so the longer corrected names are visible.
#  - Moved the active selection from K23 to B30:B33 (the corrected cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections -------------------------------------------------
# "Torra di Murtella" appears (repeated) in B18:B21.
$ws.Range("B18:B21").Value = "A Torra di Murtella"

# "Maffalcu" appears (repeated) in B30:B33.
$ws.Range("B30:B33").Value = "Malfalcu"

# --- Column width -------------------------------------------------------
# Column B grew from ~15.55 to ~32.57 (character units) to fit the longer text.
$ws.Columns.Item(2).ColumnWidth = 31.67

# --- Selection ------------------------------------------------------------
# Active selection moves onto the corrected range B30:B33 (active cell B30).
$ws.Range("B30:B33").Select()
